$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Acquire', ['{3}{U}{U}', 'Sorcery', 'Search target opponent’s library for an artifact card and put that card onto the battlefield under your control. Then that player shuffles their library.'])"
$ws.Range("A3").Value = "('Duress', ['{B}', 'Sorcery', 'Target opponent reveals their hand. You choose a noncreature, nonland card from it. That player discards that card.'])"
$ws.Range("A4").Value = "('Wash Out', ['{3}{U}', 'Sorcery', 'Return all permanents of the color of your choice to their owners’ hands.'])"

$ws.Range("A5:A13").EntireRow.Delete()
